$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.517.85'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.261.70'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '118.56'
$ws.Range("E5").Value = '  +4.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '264.33'
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("E7").Value = '  +2.30%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.36'
$ws.Range("E10").Value = '  -2.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.15'
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("E13").Value = '  -1.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.28'
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.903'
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Value = '2.602.93'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '2.267.58'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").Value = '43.500.87'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000109'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.86'
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.90'
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.69'
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.43'
$ws.Range("E24").Value = '  -5.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.86'
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("E26").Value = '  +1.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.43'
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.69'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.62'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0911'
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.68'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.129'
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("E36").Value = '  +14.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0375'
$ws.Range("E37").Value = '  +6.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.55'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.52'
$ws.Range("E40").Value = '  +4.27%  '
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.91'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.80'
$ws.Range("E42").Value = '  -7.52%  '
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -8.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '73.93'
$ws.Range("E47").Value = '  +41.35%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.49'
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.25'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0995'
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '100.71'
$ws.Range("E51").Value = '  -1.30%  '
